$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

$ws.Range("B2").Value = 10.398420022686878
$ws.Range("C2").Value = 22.384284072851585
$ws.Range("D2").Value = 28.20280377399456
$ws.Range("E2").Value = 23.381542284777765

$ws.Range("B3").Value = 9.9566180160818476
$ws.Range("C3").ClearContents()
$ws.Range("D3").Value = 38.34060839500205
$ws.Range("E3").Value = 14.184921720323047

$ws.Range("B1:E3").Select()
